# Update odds values on Sheet1 (row 2 and row 4) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 updates ---
$ws.Range("G2").Value = 1.53
$ws.Range("H2").Value = 3.8
$ws.Range("N2").Value = 8.5
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("Z2").Value = 10
$ws.Range("AC2").Value = 8.5
$ws.Range("AE2").Value = 23
$ws.Range("AT2").Value = 2.63
$ws.Range("AU2").Value = 9.5
$ws.Range("AV2").Value = 67

# --- Row 4 updates ---
$ws.Range("G4").Value = 3.45
$ws.Range("I4").Value = 2.12
$ws.Range("J4").Value = 4.1
$ws.Range("L4").Value = 2.82
$ws.Range("N4").Value = 7.3
$ws.Range("P4").Value = 2.37
$ws.Range("Q4").Value = 2.32
$ws.Range("R4").Value = 1.47
$ws.Range("S4").Value = 1.52
$ws.Range("T4").Value = 2.22
$ws.Range("U4").Value = 2.02
$ws.Range("V4").Value = 1.62
$ws.Range("X4").Value = 16.5
$ws.Range("Y4").Value = 12.5
$ws.Range("AA4").Value = 37
$ws.Range("AC4").Value = 6.8
$ws.Range("AD4").Value = 6.1
$ws.Range("AE4").Value = 18.5
$ws.Range("AF4").Value = 120
$ws.Range("AH4").Value = 5.7
$ws.Range("AI4").Value = 8.75
$ws.Range("AL4").Value = 21
$ws.Range("AN4").Value = 5.1
$ws.Range("AO4").Value = 20
$ws.Range("AR4").Value = 175
$ws.Range("AU4").Value = 8
$ws.Range("AV4").Value = 100
$ws.Range("AW4").Value = 3.8
$ws.Range("AX4").Value = 11.5
$ws.Range("AY4").Value = 25
$ws.Range("BA4").Value = 110
$ws.Range("BB4").Value = 450
